# Applies the "corrected BidExclusion apply-to-all issue" edit:
#  - Results sheet: Bid ID 2 (row 3) is no longer excluded (gets real award values);
#    Bid IDs 1,3,4,5,6,7,8,9,10 (rows 2,4-11) become "No Bid" / zeroed out.
#  - Feasibility Notes sheet: updated message explaining infeasibility.
#  - LP Model sheet: BidExclusion_0_* constraints corrected to only exclude Bid ID 2
#    (for suppliers A, B and C) instead of excluding every Bid ID for suppliers B/C.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Results sheet
# ---------------------------------------------------------------------------
$results = $wb.Worksheets.Item("Results")

# Column I ("Percentage Volume Discount") holds text values like "1%"/"0%".
# Force the column to a text format first so assigning strings such as "0%"
# is not auto-parsed by Excel into a numeric percentage.
$results.Range("I2:I11").NumberFormat = "@"

# Rows that become "No Bid" (award pulled): row 2 (Bid ID 1) and rows 4-11 (Bid IDs 3-10)
$noBidRows = @(2,4,5,6,7,8,9,10,11)
foreach ($r in $noBidRows) {
    $results.Cells.Item($r, 6).Value2  = 0        # F: Baseline Spend
    $results.Cells.Item($r, 7).Value2  = "No Bid" # G: Awarded Supplier
    $results.Cells.Item($r, 8).Value2  = 0        # H: Original Awarded Supplier Price
    $results.Cells.Item($r, 9).Value2  = "0%"     # I: Percentage Volume Discount
    $results.Cells.Item($r, 10).Value2 = 0        # J: Discounted Awarded Supplier Price
    $results.Cells.Item($r, 11).Value2 = 0        # K: Awarded Supplier Spend
    $results.Cells.Item($r, 12).Value2 = 0        # L: Awarded Volume
    $results.Cells.Item($r, 13).Value2 = 0        # M: Baseline Savings
}

# Row 3 (Bid ID 2) is no longer excluded -> full undiscounted award values
$results.Cells.Item(3, 9).Value2  = "0%"     # I3: Percentage Volume Discount
$results.Cells.Item(3, 10).Value2 = 70       # J3: Discounted Awarded Supplier Price
$results.Cells.Item(3, 11).Value2 = 630000   # K3: Awarded Supplier Spend
$results.Cells.Item(3, 13).Value2 = 774000   # M3: Baseline Savings

# ---------------------------------------------------------------------------
# 2. Feasibility Notes sheet
# ---------------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Feasibility Notes")
$nl = [char]10
$feasibilityMessage = "Model is infeasible. Likely causes include:" + $nl + `
    " - Insufficient supplier capacity relative to demand." + $nl + `
    " - Custom rule constraints conflicting with overall volume/demand." + $nl + `
    "Detailed Rule Evaluations:" + $nl + `
    "Rule 1 ('Exclude Bids'): The exclusion criteria might be too broad, removing all valid bids needed to satisfy demand for a Bid ID or grouping." + $nl + `
    "Please review supplier capacities, demand figures, and custom rule constraints for adjustments."
$notes.Range("A2").Value2 = $feasibilityMessage

# ---------------------------------------------------------------------------
# 3. LP Model sheet - fix BidExclusion constraints
# ---------------------------------------------------------------------------
$lp = $wb.Worksheets.Item("LP Model")
$lpText = $lp.Range("A2").Value2

$oldBlock = "BidExclusion_0_10_B: x_B_10 = 0" + $nl + `
    "BidExclusion_0_10_C: x_C_10 = 0" + $nl + `
    "BidExclusion_0_2_B: x_B_2 = 0" + $nl + `
    "BidExclusion_0_2_C: x_C_2 = 0" + $nl + `
    "BidExclusion_0_3_B: x_B_3 = 0" + $nl + `
    "BidExclusion_0_3_C: x_C_3 = 0" + $nl + `
    "BidExclusion_0_4_B: x_B_4 = 0" + $nl + `
    "BidExclusion_0_4_C: x_C_4 = 0" + $nl + `
    "BidExclusion_0_5_B: x_B_5 = 0" + $nl + `
    "BidExclusion_0_5_C: x_C_5 = 0" + $nl + `
    "BidExclusion_0_6_B: x_B_6 = 0" + $nl + `
    "BidExclusion_0_6_C: x_C_6 = 0" + $nl + `
    "BidExclusion_0_7_B: x_B_7 = 0" + $nl + `
    "BidExclusion_0_7_C: x_C_7 = 0" + $nl + `
    "BidExclusion_0_8_B: x_B_8 = 0" + $nl + `
    "BidExclusion_0_8_C: x_C_8 = 0" + $nl + `
    "BidExclusion_0_9_B: x_B_9 = 0" + $nl + `
    "BidExclusion_0_9_C: x_C_9 = 0" + $nl

$newBlock = "BidExclusion_0_2_A: x_A_2 = 0" + $nl + `
    "BidExclusion_0_2_B: x_B_2 = 0" + $nl + `
    "BidExclusion_0_2_C: x_C_2 = 0" + $nl

if ($lpText.IndexOf($oldBlock) -lt 0) {
    throw "BidExclusion block not found in LP Model text; cannot apply fix."
}
$lpText = $lpText.Replace($oldBlock, $newBlock)
$lp.Range("A2").Value2 = $lpText
